# Agenda slide ("RXJS Make it easy"): capitalize the four RxJS mapping
# operators listed in the "Common operators in RXJS" bullet and turn the
# single run into individually-capitalized, comma-separated runs, e.g.
#   switchMap, exhaustMap, mergeMap, concatMap
# becomes
#   SwitchMap, ExhaustMap, MergeMap, ConcatMap
#
# The slide is slide #2 (1-based), the shape is the body "Content
# Placeholder" (shape #2). We locate each lowercase operator name inside
# the shape's running text and overwrite just that word in place via
# TextRange.Characters(start, length) — this mirrors how PowerPoint
# itself turns a single run into several runs when only part of its text
# is retyped, while leaving every other run/paragraph untouched.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(2)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

function Set-WordText($range, [string]$oldWord, [string]$newWord) {
    $current = $range.Text
    $idx = $current.IndexOf($oldWord)
    if ($idx -lt 0) {
        throw "Could not find '$oldWord' in shape text"
    }
    $wordRange = $range.Characters($idx + 1, $oldWord.Length)
    $wordRange.Text = $newWord
}

Set-WordText $textRange "switchMap" "SwitchMap"
Set-WordText $textRange "exhaustMap" "ExhaustMap"
Set-WordText $textRange "mergeMap" "MergeMap"
Set-WordText $textRange "concatMap" "ConcatMap"
